$d = $word.ActiveDocument

# The document currently ends with an empty list paragraph (ilvl 0) whose
# only content is the _GoBack bookmark. That paragraph becomes the first of
# three new bold list items; two more new list items are appended after it,
# the last one nested at ilvl 1. The _GoBack bookmark ends up inside that
# last new paragraph, positioned between its two runs.

# Remove the existing _GoBack bookmark; it will be re-created once the new
# paragraph text exists, in its final resting place.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# Paragraph 3 used to hold only the bookmark -- give it its text (bold).
$p3 = $d.Paragraphs(3)
$p3.Range.InsertAfter("How to list existing variables for user to select from when calculating or assigning?")
$p3.Range.Font.Bold = 1

# Insert a new list paragraph after it (still ilvl 0), bold text.
$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs(4)
$p4.Range.InsertAfter("Ability to assign a “statement” as an entry to another statement (nested statements)")
$p4.Range.Font.Bold = 1

# Insert a further new list paragraph after that one, nested at ilvl 1. Put
# both runs of text in now (with the trailing "?") so that the bookmark we
# add next lands strictly between the two runs instead of right on the
# paragraph mark.
$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs(5)
$p5.Range.InsertAfter("Have popup box which has our normal statement buttons at the top?")
$p5.Range.Font.Bold = 1
$p5.Range.ListFormat.ListLevelNumber = 2

# Re-create the _GoBack bookmark collapsed right before the trailing "?"
# (i.e. immediately after "...top").
$p5 = $d.Paragraphs(5)
$bmPos = $p5.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
